$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2530.2856
$ws.Range("J4").Value = 3454.8
$ws.Range("L4").Value = 3454.8
$ws.Range("N4").Value = -3682.8
$ws.Range("H38").Value = 1191.1666
$ws.Range("I38").Value = 679.4
$ws.Range("J38").Value = 3750
$ws.Range("K38").Value = 2038.2
$ws.Range("L38").Value = 11250
$ws.Range("M38").Value = -1666.2
$ws.Range("N38").Value = -11994
$ws.Range("H58").Value = 6386.125
$ws.Range("J58").Value = 25000
$ws.Range("L58").Value = 75000
$ws.Range("N58").Value = -75300
$ws.Range("H87").Value = 64746.668
$ws.Range("J87").Value = 64746.668
$ws.Range("L87").Value = 64746.668
$ws.Range("N87").Value = -67242.66800000001
$ws.Range("H90").Value = 64746.668
$ws.Range("J90").Value = 64746.668
$ws.Range("L90").Value = 194240.004
$ws.Range("N90").Value = -206720.004
$ws.Range("H106").Value = 15489.444
$ws.Range("I106").Value = 3135
$ws.Range("K106").Value = 3135
$ws.Range("M106").Value = -2504
$ws.Range("H115").Value = 655.6667
$ws.Range("I115").Value = 655.6667
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1967.0001
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -400.0001
$ws.Range("N115").Value = ""
$ws.Range("H127").Value = 672.5
$ws.Range("I127").Value = 672.5
$ws.Range("K127").Value = 2017.5
$ws.Range("M127").Value = 2942.5
$ws.Range("H137").Value = 5766
$ws.Range("I137").Value = 1725.25
$ws.Range("K137").Value = 5175.75
$ws.Range("M137").Value = -2625.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6671615
$ws.Range("I32").Value = 8335355.5
$ws.Range("K32").Value = 8335355.5
$ws.Range("M32").Value = -8335068.5
$ws.Range("H61").Value = 39484604
$ws.Range("I61").Value = 41677170
$ws.Range("J61").Value = 35725930
$ws.Range("K61").Value = 41677170
$ws.Range("L61").Value = 35725930
$ws.Range("M61").Value = -41676958
$ws.Range("N61").Value = -35726354
$ws.Range("H117").Value = 104990
$ws.Range("J117").Value = 104990
$ws.Range("L117").Value = 104990
$ws.Range("N117").Value = -114168
$ws.Range("H122").Value = 2663.4644
$ws.Range("I122").Value = 1772.4736
$ws.Range("K122").Value = 5317.4208
$ws.Range("M122").Value = -2867.4208
$ws.Range("H136").Value = 39484604
$ws.Range("I136").Value = 41677170
$ws.Range("J136").Value = 35725930
$ws.Range("K136").Value = 125031510
$ws.Range("L136").Value = 107177790
$ws.Range("M136").Value = -125028960
$ws.Range("N136").Value = -107182890

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 480805.03
$ws.Range("I134").Value = 1796.2354
$ws.Range("J134").Value = 1837996.6
$ws.Range("K134").Value = 5388.706200000001
$ws.Range("L134").Value = 5513989.800000001
$ws.Range("M134").Value = -2853.706200000001
$ws.Range("N134").Value = -5519059.800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 158.8
$ws.Range("I22").Value = 158.8
$ws.Range("K22").Value = 158.8
$ws.Range("M22").Value = 191.2
$ws.Range("H141").Value = 334833.2
$ws.Range("J141").Value = 358416.5
$ws.Range("L141").Value = 358416.5
$ws.Range("N141").Value = -368776.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 145.08333
$ws.Range("I2").Value = 91.333336
$ws.Range("K2").Value = 548.000016
$ws.Range("M2").Value = -435.000016
$ws.Range("H17").Value = 440
$ws.Range("I17").Value = 200
$ws.Range("K17").Value = 600
$ws.Range("M17").Value = -431
$ws.Range("H34").Value = 381
$ws.Range("I34").Value = 381
$ws.Range("K34").Value = 1143
$ws.Range("M34").Value = -1059
$ws.Range("H39").Value = 2332.8333
$ws.Range("J39").Value = 3165.6667
$ws.Range("L39").Value = 9497.000100000001
$ws.Range("N39").Value = -10085.0001
$ws.Range("H55").Value = 2320.6
$ws.Range("J55").Value = 1500
$ws.Range("L55").Value = 4500
$ws.Range("N55").Value = -4854
$ws.Range("H103").Value = 1623.1
$ws.Range("I103").Value = 395.25
$ws.Range("J103").Value = 2441.6667
$ws.Range("K103").Value = 1185.75
$ws.Range("L103").Value = 7325.000100000001
$ws.Range("M103").Value = -306.75
$ws.Range("N103").Value = -9083.000100000001
$ws.Range("H117").Value = 1395.1428
$ws.Range("I117").Value = 850
$ws.Range("J117").Value = 1486
$ws.Range("K117").Value = 2550
$ws.Range("L117").Value = 4458
$ws.Range("M117").Value = 892
$ws.Range("N117").Value = -11342

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 2349.5
$ws.Range("J17").Value = 3066
$ws.Range("L17").Value = 3066
$ws.Range("N17").Value = -3402
$ws.Range("H62").Value = 1011000
$ws.Range("J62").Value = 1011000
$ws.Range("L62").Value = 1011000
$ws.Range("N62").Value = -1012372
$ws.Range("H65").Value = 1011000
$ws.Range("J65").Value = 1011000
$ws.Range("L65").Value = 3033000
$ws.Range("N65").Value = -3039864
$ws.Range("H80").Value = 5216.1113
$ws.Range("I80").Value = 4667.6665
$ws.Range("J80").Value = 5490.3335
$ws.Range("K80").Value = 4667.6665
$ws.Range("L80").Value = 5490.3335
$ws.Range("M80").Value = -3669.6665
$ws.Range("N80").Value = -7486.3335
$ws.Range("H83").Value = 5216.1113
$ws.Range("I83").Value = 4667.6665
$ws.Range("J83").Value = 5490.3335
$ws.Range("K83").Value = 23338.3325
$ws.Range("L83").Value = 27451.6675
$ws.Range("M83").Value = -18346.3325
$ws.Range("N83").Value = -37435.6675
$ws.Range("H102").Value = 3011.75
$ws.Range("I102").Value = 2521.375
$ws.Range("J102").Value = 3665.5833
$ws.Range("K102").Value = 2521.375
$ws.Range("L102").Value = 3665.5833
$ws.Range("M102").Value = -899.375
$ws.Range("N102").Value = -6909.5833
$ws.Range("H122").Value = 2101.8667
$ws.Range("I122").Value = 2149.1667
$ws.Range("J122").Value = 1912.6666
$ws.Range("K122").Value = 6447.500100000001
$ws.Range("L122").Value = 5737.9998
$ws.Range("M122").Value = -3997.500100000001
$ws.Range("N122").Value = -10637.9998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1219.6666
$ws.Range("I16").Value = 1485.75
$ws.Range("K16").Value = 1485.75
$ws.Range("M16").Value = -1315.75
$ws.Range("H35").Value = 9500
$ws.Range("I35").Value = 9500
$ws.Range("K35").Value = 9500
$ws.Range("M35").Value = -9164
$ws.Range("H122").Value = 6820.763
$ws.Range("I122").Value = 6686.7856
$ws.Range("K122").Value = 20060.3568
$ws.Range("M122").Value = -17610.3568
$ws.Range("H136").Value = 66042.62
$ws.Range("I136").Value = 9271.666999999999
$ws.Range("K136").Value = 27815.001
$ws.Range("M136").Value = -25265.001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 27897.8
$ws.Range("I54").Value = 17250
$ws.Range("J54").Value = 34996.332
$ws.Range("K54").Value = 17250
$ws.Range("L54").Value = 34996.332
$ws.Range("M54").Value = -16730
$ws.Range("N54").Value = -36036.332
$ws.Range("H126").Value = 3194.5715
$ws.Range("I126").Value = 3059.5
$ws.Range("J126").Value = 4005
$ws.Range("K126").Value = 9178.5
$ws.Range("L126").Value = 12015
$ws.Range("M126").Value = -6708.5
$ws.Range("N126").Value = -16955
$ws.Range("H136").Value = 5109.2593
$ws.Range("I136").Value = 5365.905
$ws.Range("J136").Value = 4211
$ws.Range("K136").Value = 16097.715
$ws.Range("L136").Value = 12633
$ws.Range("M136").Value = -13547.715
$ws.Range("N136").Value = -17733

Write-Host "Applied all changes"